$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (BuiltInSupport.c) - Reviewed / RLB / No changes required.
$ws.Range("B12").Value = "Reviewed"
$ws.Range("C12").Value = "RLB"
$ws.Range("D12").Value = "No changes required."

# Row 21 (HLUFunctions.c) - Patched / RLB / new note about variable "total"
$ws.Range("B21").Value = "Patched"
$ws.Range("C21").Value = "RLB"
$ws.Range("D21").Value = "Many changes, mostly to promote variable " + [char]8220 + "total" + [char]8221 + " throughout. _NclAddToOverlay2() severely dysfunctional!"

# Row 22 (HLUFunctions.h)
$ws.Range("B22").Value = "Reviewed"
$ws.Range("C22").Value = "RLB"
$ws.Range("D22").Value = "No changes required."

# Row 23 (HLUSupport.c)
$ws.Range("B23").Value = "Reviewed"
$ws.Range("C23").Value = "RLB"
$ws.Range("D23").Value = "No changes required."

# Row 24 (HLUSupport.h)
$ws.Range("B24").Value = "Reviewed"
$ws.Range("C24").Value = "RLB"
$ws.Range("D24").Value = "No changes required."
